$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.434.68'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '3.099.15'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.48'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.88%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '3.088.41'
$ws.Range("E8").Value = '  -0.72%  '

$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.75%  '

$ws.Range("E11").Value = '  -2.68%  '

$ws.Range("E12").Value = '  -2.79%  '

$ws.Range("E13").Value = '  -1.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.34'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.79%  '

$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("D16").Value = '3.611.49'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").Value = '63.301.43'
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.08'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("D19").Value = '3.092.62'
$ws.Range("E19").Value = '  -0.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '460.82'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.21'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.80%  '

$ws.Range("E22").Value = '  -0.74%  '

$ws.Range("E23").Value = '  -1.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.98'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +8.42%  '

$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.66'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.20'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.81'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.110'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.61'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.66%  '

$ws.Range("D35").Value = '0.0₃0849'
$ws.Range("E35").Value = '  -3.12%  '

$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("E37").Value = '  -3.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.34'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.08%  '

$ws.Range("E39").Value = '  -1.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.13'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '434.03'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.72'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0367'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.21%  '

$ws.Range("D44").Value = '2.870.53'
$ws.Range("E44").Value = '  -2.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.270'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.25%  '

$ws.Range("E46").Value = '  -2.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.96'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.67'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("E50").Value = '  -1.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.04'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.65%  '
